$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values
$ws.Range("B2").Value = 5.2558677754913283
$ws.Range("C2").Value = 5.9266320298577853
$ws.Range("D2").Value = 3.1094321911993568
$ws.Range("E2").Value = 3.3615954877177785

# Update row 3 values
$ws.Range("B3").Value = 4.678302885769237
$ws.Range("C3").Value = 8.1776421550476446
$ws.Range("D3").Value = 2.8182186993434706
$ws.Range("E3").Value = 0.71005297835652481

# Update selection to reflect the new selected range
$ws.Range("B1:E3").Select()
